$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the last three data rows (former rows 8, 9, 10 - "ARUN", "ANANDHU", "FAHAD")
$ws.Rows("8:10").Delete()

# Row 3 - SHABIN
$ws.Range("B3").Value = "16-12-2025"
$ws.Range("C3").Value = "SHABIN"
$ws.Range("D3").Value = 8129192047
$ws.Range("E3").Value = "18-01-2026"
$ws.Range("F3").Value = "ARJUN G.S"
$ws.Range("H3").Value = "ENQUIRY"
$ws.Range("I3").Value = "ENQUIRY WITHOUT TRIAL"
$ws.Range("K3").Value = "WILL CONFIRM TOMO"

# Row 4 - JAKSON
$ws.Range("B4").Value = "16-12-2025"
$ws.Range("C4").Value = "JAKSON"
$ws.Range("D4").Value = 9567658570
$ws.Range("E4").Value = "29-12-2025"
$ws.Range("F4").Value = "NIHAL S"
$ws.Range("I4").Value = "ENQUIRY WITHOUT TRIAL"
$ws.Range("K4").Value = "WILL CONFIRM LATER"

# Row 5 - ABI
$ws.Range("B5").Value = "16-12-2025"
$ws.Range("C5").Value = "ABI"
$ws.Range("D5").Value = 7907547616
$ws.Range("E5").Value = "24-01-2026"
$ws.Range("F5").Value = "MOHAMMED NABEEL N"
$ws.Range("H5").Value = "ENQUIRY"
$ws.Range("I5").Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Range("K5").ClearContents()

# Row 6 - ABISHEK
$ws.Range("B6").Value = "16-12-2025"
$ws.Range("C6").Value = "ABISHEK"
$ws.Range("D6").Value = 9656112979
# "04-01-2026" would otherwise be auto-recognized as a date (MM-DD-YYYY); force plain text
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "04-01-2026"
$ws.Range("E6").Style = "Normal"
$ws.Range("F6").Value = "ARJUN G.S"
$ws.Range("H6").Value = "ENQUIRY"
$ws.Range("I6").Value = "ENQUIRY WITHOUT TRIAL"
$ws.Range("K6").ClearContents()

# Row 7 - NASIM
$ws.Range("B7").Value = "16-12-2025"
$ws.Range("C7").Value = "NASIM"
$ws.Range("D7").Value = 8113969331
# "12-01-2026" would otherwise be auto-recognized as a date (MM-DD-YYYY); force plain text
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "12-01-2026"
$ws.Range("E7").Style = "Normal"
$ws.Range("F7").Value = "NIHAL S"
$ws.Range("H7").Value = "ENQUIRY"
$ws.Range("I7").Value = "ENQUIRY WITHOUT BRIDE/FAMILY"
$ws.Range("K7").Value = "WILL VISIT LATER"

# Column width adjustments
$ws.Columns(8).ColumnWidth = 10.8
$ws.Columns(11).ColumnWidth = 24.3
